$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster table (Oyuncu Adı / Pozisyon / Takım)
$data = @(
    @("Norman Powell",     "SG,SF",      "LA Clippers"),
    @("Max Christie",      "SG,SF",      "Dallas Mavericks"),
    @("Jimmy Butler",      "SF,PF",      "Golden State Warriors"),
    @("Nikola Jovic",      "PF,C",       "Miami Heat"),
    @("Kel'el Ware",       "PF,C",       "Miami Heat"),
    @("Walker Kessler",    "C",          "Utah Jazz"),
    @("Jalen Williams",    "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Jalen Brunson",     "PG",         "New York Knicks"),
    @("Trae Young",        "PG",         "Atlanta Hawks"),
    @("Bol Bol",           "PF,C",       "Phoenix Suns"),
    @("Christian Braun",   "SG,SF",      "Denver Nuggets"),
    @("LeBron James",      "SF,PF",      "Los Angeles Lakers"),
    @("Kawhi Leonard",     "SG,SF,PF",   "LA Clippers"),
    @("Desmond Bane",      "SG,SF",      "Memphis Grizzlies"),
    @("Immanuel Quickley", "PG,SG",      "Toronto Raptors"),
    @("Myles Turner",      "C",          "Indiana Pacers"),
    @("Brandon Ingram",    "PG,SG",      "Toronto Raptors"),
    @("Devin Booker",      "PG,SG",      "Phoenix Suns")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row = $row + 1
}
